# Weekly update: insert 3 new price rows (week of 2023-02-07, serial 44964)
# right before row 351, shifting the existing rows 351:453 down to 354:456.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 351, pushing 351:453 down to 354:456.
$ws.Rows("351:353").Insert()

# Columns A,B,C,E,F,G,H,I,J,K,T are constant across every data row in this sheet.
$commonA = 7
$commonB = "Terminal Hortofrutícola Agro Chillán"
$commonC = "Ñuble"
$commonE = 16
$commonF = "Fruta"
$commonG = 100101
$commonH = "Berries"
$commonI = 100112025
$commonJ = "Frutilla"
$commonK = "Sin especificar"
$commonT = 7

# New rows' per-row data (Fecha serial, Calidad, Volumen, Precio min/max/promedio,
# Unidad, Origen, Precio $/Kg).
$newRows = @(
    @{ Row = 351; D = 44964; L = "Especial"; M = 50; N = 7500; O = 7500; P = 7500; Q = "`$/caja 7 kilos"; R = "Provincia de Diguillín"; S = 1071 },
    @{ Row = 352; D = 44964; L = "Primera";  M = 30; N = 6500; O = 6500; P = 6500; Q = "`$/caja 7 kilos"; R = "Provincia de Diguillín"; S = 929 },
    @{ Row = 353; D = 44964; L = "Segunda";  M = 80; N = 5000; O = 5500; P = 5188; Q = "`$/caja 7 kilos"; R = "Provincia de Diguillín"; S = 741 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $commonA
    $ws.Cells.Item($row, 2).Value = $commonB
    $ws.Cells.Item($row, 3).Value = $commonC
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $commonE
    $ws.Cells.Item($row, 6).Value = $commonF
    $ws.Cells.Item($row, 7).Value = $commonG
    $ws.Cells.Item($row, 8).Value = $commonH
    $ws.Cells.Item($row, 9).Value = $commonI
    $ws.Cells.Item($row, 10).Value = $commonJ
    $ws.Cells.Item($row, 11).Value = $commonK
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $commonT
}
